# 20231224 Revision_Simulated_Tornado_15M.xlsx
# "50 simulation for 0M budget"
#
# The "15M-base-worst" sheet picks a random Block ID via K11 = RANDBETWEEN(1,100)
# and propagates it through L11/M11 (lookup), L12/M12 (scaling) and the O16:O115
# distance column. The target commit re-rolled that random draw so that K11 lands
# on Block ID 1 (all the downstream cells are pure recalculations of the same
# formulas against that new K11).
#
# RANDBETWEEN is volatile, so simply writing "1" into K11 would either (a) blow
# away the formula if we use .Value, or (b) get re-rolled again by the automatic
# recalculation that happens once this script returns. Instead we keep
# re-entering the live formula until it happens to settle on the value we need,
# then flip calculation to manual so the trailing auto-recalc leaves the cached
# results (K11 and everything downstream of it) alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("15M-base-worst")

$target = "1"
$iterations = 0
$current = ""
do {
    $ws.Range("K11").Formula = "=RANDBETWEEN(1,100)"
    $current = $ws.Range("K11").Text
    $iterations = $iterations + 1
} while ($current -ne $target -and $iterations -lt 100000)

# Freeze calculation so the post-script recalc (which would otherwise re-roll
# the volatile RANDBETWEEN) leaves the cells we just computed untouched.
$excel.Calculation = -4135

# Results sheet: the active selection moved from E53 to F53 (and Excel dropped
# the saved scroll position topLeftCell="A52" along with it).
$rs = $wb.Worksheets.Item("Results")
$rs.Range("F53").Select()
